$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the 380 kV case (res_line/pl_mw), rows 2-25,
# columns C,D,E,F,G,I,K,M. Columns B,H,J,L,N,O remain 0 and are untouched.
$data = @{
    2 = @{ "C"=0.01644600113095862; "D"=0.05921444607283988; "E"=0.07458911246050093; "F"=1.734368667262302; "G"=0.002477133199180375; "I"=1.335212818010618; "K"=1.352155494448709; "M"=0.4203527620946375 }
    3 = @{ "C"=0.01643156202132268; "D"=0.05961461513903288; "E"=0.06935873065619447; "F"=1.686935561353067; "G"=0.002482597736616221; "I"=1.2994965734689; "K"=1.221517008179603; "M"=0.3828091206204647 }
    4 = @{ "C"=0.01642105954081252; "D"=0.05988410839940883; "E"=0.06620178186767589; "F"=1.658966743776446; "G"=0.002486125618293021; "I"=1.278456008639736; "K"=1.141922323139624; "M"=0.3599884099054194 }
    5 = @{ "C"=0.01641637541846208; "D"=0.05999984943760239; "E"=0.06492871164675407; "F"=1.647856807896815; "G"=0.002487606826465839; "I"=1.270103165308853; "K"=1.109640783929819; "M"=0.350746087944195 }
    6 = @{ "C"=0.01641557342433941; "D"=0.06001942393595883; "E"=0.06471812065857918; "F"=1.646029295635486; "G"=0.002487855416031205; "I"=1.268729485113298; "K"=1.104289716205017; "M"=0.3492148440783467 }
    7 = @{ "C"=0.01642099799625818; "D"=0.0598856454342318; "E"=0.06618455888439456; "F"=1.658815750760823; "G"=0.002486145417951358; "I"=1.278342466026999; "K"=1.141486341551058; "M"=0.3598635338004073 }
    8 = @{ "C"=0.01644136656754291; "D"=0.05934745027062149; "E"=0.0727741775300359; "F"=1.717772167615323; "G"=0.002478981636499864; "I"=1.322711884735114; "K"=1.306982185202003; "M"=0.407359196991834 }
    9 = @{ "C"=0.01646800784610036; "D"=0.05848354889464957; "E"=0.08614325135629741; "F"=1.84269189390352; "G"=0.002466296045941101; "I"=1.416885275314257; "K"=1.63650626254389; "M"=0.5023764633544658 }
    10 = @{ "C"=0.01647906014092548; "D"=0.0579693752878363; "E"=0.09625862483905934; "F"=1.940348117315978; "G"=0.00245779647300577; "I"=1.490601006326074; "K"=1.881794495536042; "M"=0.5733993375390298 }
    11 = @{ "C"=0.01648215450957657; "D"=0.05776246223371118; "E"=0.1009285874198298; "F"=1.986095767998677; "G"=0.002454105798873353; "I"=1.525154009565668; "K"=1.994108910505361; "M"=0.6059885930825573 }
    12 = @{ "C"=0.01648304174253745; "D"=0.05768805943824162; "E"=0.102707176859532; "F"=2.003612980429864; "G"=0.002452733353961105; "I"=1.538387611226369; "K"=2.0367469976282; "M"=0.6183707836078156 }
    13 = @{ "C"=0.01648286341707816; "D"=0.05770390658214808; "E"=0.1023236678264396; "F"=1.999831673737418; "G"=0.002453027818980078; "I"=1.535530844916977; "K"=2.027559330935844; "M"=0.6157022022880909 }
    14 = @{ "C"=0.01648223324098197; "D"=0.05775626149058155; "E"=0.1010747071617288; "F"=1.987533021305921; "G"=0.002453992384279072; "I"=1.526239742481806; "K"=1.997614615120995; "M"=0.6070064488994404 }
    15 = @{ "C"=0.01648180999586657; "D"=0.05778884696871955; "E"=0.1003110179487976; "F"=1.980025047836818; "G"=0.002454586476739364; "I"=1.520568171049021; "K"=1.979286610318752; "M"=0.6016854634449516 }
    16 = @{ "C"=0.01647881844983345; "D"=0.0579834465570741; "E"=0.09595483738580413; "F"=1.93738530902553; "G"=0.002458041192438497; "I"=1.488363615452073; "K"=1.874469362371372; "M"=0.5712752840594106 }
    17 = @{ "C"=0.01647648333672436; "D"=0.05810978869886796; "E"=0.09330022563766249; "F"=1.911568470553192; "G"=0.002460205474117329; "I"=1.468870048416846; "K"=1.810356044630794; "M"=0.5526921932475091 }
    18 = @{ "C"=0.01647495863885773; "D"=0.05818499328385229; "E"=0.0917797836345855; "F"=1.896843669702747; "G"=0.002461466870149513; "I"=1.457753656238054; "K"=1.773548502727465; "M"=0.5420300406400003 }
    19 = @{ "C"=0.01647441141330752; "D"=0.05821088969173616; "E"=0.09126608064881481; "F"=1.891879365133065; "G"=0.002461896805508419; "I"=1.45400621196417; "K"=1.761097852796013; "M"=0.5384245145465627 }
    20 = @{ "C"=0.01647675075694544; "D"=0.05809607639892889; "E"=0.0935821466757929; "F"=1.914303822491689; "G"=0.002459973369829171; "I"=1.470935239420371; "K"=1.817173888936679; "M"=0.5546676607398098 }
    21 = @{ "C"=0.0164824261126153; "D"=0.0577407758032038; "E"=0.1014412783125351; "F"=1.991140153955797; "G"=0.002453708387257315; "I"=1.528964696728039; "K"=2.00640718337678; "M"=0.6095594722512203 }
    22 = @{ "C"=0.01648447392021346; "D"=0.05753162695759073; "E"=0.1066371474838945; "F"=2.042486955369128; "G"=0.002449760283312218; "I"=1.567760741629542; "K"=2.13070700087826; "M"=0.6456759246390362 }
    23 = @{ "C"=0.01648353507624911; "D"=0.05764112026702506; "E"=0.1038584604497785; "F"=2.014977732270296; "G"=0.002451854112352604; "I"=1.546974068434722; "K"=2.064308034705277; "M"=0.6263774629913712 }
    24 = @{ "C"=0.01647663042393255; "D"=0.05810226773444782; "E"=0.0934546723530616; "F"=1.913066803694306; "G"=0.002460078250727778; "I"=1.470001284190076; "K"=1.814091379514537; "M"=0.5537744849507504 }
    25 = @{ "C"=0.01646226816914975; "D"=0.05869636629100938; "E"=0.08247659784863259; "F"=1.807880471343026; "G"=0.002469583011092047; "I"=1.390625498502203; "K"=1.546812476570722; "M"=0.4764636719579016 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
